# "Añadida la logica hasta 2020"
# Extends the MOTOGP file list with the 2019 and 2020 entries and syncs
# Hoja1 (the "active" single-row lookup sheet) with the full Hoja2 list.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# ---------------------------------------------------------------------
# 1) Hoja2: append the 2018, 2019 and 2020 rows (2002-2017 already exist)
# ---------------------------------------------------------------------
$ws2.Range("B18:B20").NumberFormat = "@"

$ws2.Range("A18").Value2 = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2018.xlsx"
$ws2.Range("B18").Value2 = "2018"

$ws2.Range("A19").Value2 = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2019.xlsx"
$ws2.Range("B19").Value2 = "2019"

$ws2.Range("B20").Value2 = "2020"
$ws2.Range("A20").Value2 = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2020.xlsx"

# ---------------------------------------------------------------------
# 2) Hoja1: rebuild rows 2-20 so it mirrors Hoja2's full 2002-2020 list
# ---------------------------------------------------------------------
$years = 2002..2020

$ws1.Range("B2:B20").NumberFormat = "@"

for ($i = 0; $i -lt $years.Count; $i++) {
    $year = $years[$i]
    $row  = $i + 2
    $ws1.Range("A$row").Value2 = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\$year.xlsx"
    $ws1.Range("B$row").Value2 = "$year"
}

# ---------------------------------------------------------------------
# 3) Selections, matching the saved state in the diff
# ---------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("A2:B20").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("A24").Select() | Out-Null
